$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row 68: date, hours, running total
$ws.Range("A68").Value = 45443
$ws.Range("A67").Copy()
$ws.Range("A68").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("B68").Value = 3.5

$ws.Range("C68").Formula = "=C67+B68"

# Update selection/view to match new last row
$ws.Range("C68").Select()

$wb.Save()
